{"js": "// Replace the date line and each division problem's text with its new\n// value. Every old string in the mapping is unique within the document,\n// so a plain body.search + insertText(replace) round-trip is sufficient\n// and preserves each run's existing formatting (fonts, size, etc.).\nconst replacements = [\n  [\"2024-11-12 Tuesday\", \"2024-11-13 Wednesday\"],\n  [\"10\u00f77=\", \"73\u00f78=\"],\n  [\"13\u00f77=\", \"23\u00f78=\"],\n  [\"99\u00f75=\", \"10\u00f75=\"],\n  [\"11\u00f79=\", \"31\u00f73=\"],\n  [\"79\u00f72=\", \"67\u00f77=\"],\n  [\"33\u00f79=\", \"30\u00f73=\"],\n  [\"18\u00f78=\", \"28\u00f72=\"],\n  [\"26\u00f75=\", \"65\u00f77=\"],\n  [\"63\u00f74=\", \"89\u00f77=\"],\n  [\"71\u00f79=\", \"29\u00f75=\"],\n  [\"51\u00f75=\", \"29\u00f73=\"],\n  [\"52\u00f75=\", \"90\u00f78=\"],\n  [\"16\u00f77=\", \"43\u00f78=\"],\n  [\"71\u00f73=\", \"77\u00f77=\"],\n  [\"56\u00f75=\", \"98\u00f78=\"],\n  [\"13\u00f79=\", \"97\u00f77=\"],\n  [\"94\u00f76=\", \"74\u00f75=\"],\n  [\"27\u00f72=\", \"18\u00f73=\"],\n  [\"61\u00f77=\", \"68\u00f75=\"],\n  [\"30\u00f74=\", \"10\u00f76=\"],\n  [\"32\u00f73=\", \"44\u00f75=\"],\n  [\"56\u00f79=\", \"41\u00f72=\"],\n  [\"89\u00f78=\", \"74\u00f78=\"],\n  [\"26\u00f79=\", \"81\u00f78=\"],\n  [\"60\u00f76=\", \"48\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each division problem's text with its new\n# value via Find & Replace. Every old string below is unique within the\n# document, so a simple Find.Execute(Replace:=wdReplaceAll) per pair is\n# sufficient and leaves each run's existing formatting (fonts, size, etc.)\n# untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-11-12 Tuesday\", \"2024-11-13 Wednesday\"),\n  @(\"10\u00f77=\", \"73\u00f78=\"),\n  @(\"13\u00f77=\", \"23\u00f78=\"),\n  @(\"99\u00f75=\", \"10\u00f75=\"),\n  @(\"11\u00f79=\", \"31\u00f73=\"),\n  @(\"79\u00f72=\", \"67\u00f77=\"),\n  @(\"33\u00f79=\", \"30\u00f73=\"),\n  @(\"18\u00f78=\", \"28\u00f72=\"),\n  @(\"26\u00f75=\", \"65\u00f77=\"),\n  @(\"63\u00f74=\", \"89\u00f77=\"),\n  @(\"71\u00f79=\", \"29\u00f75=\"),\n  @(\"51\u00f75=\", \"29\u00f73=\"),\n  @(\"52\u00f75=\", \"90\u00f78=\"),\n  @(\"16\u00f77=\", \"43\u00f78=\"),\n  @(\"71\u00f73=\", \"77\u00f77=\"),\n  @(\"56\u00f75=\", \"98\u00f78=\"),\n  @(\"13\u00f79=\", \"97\u00f77=\"),\n  @(\"94\u00f76=\", \"74\u00f75=\"),\n  @(\"27\u00f72=\", \"18\u00f73=\"),\n  @(\"61\u00f77=\", \"68\u00f75=\"),\n  @(\"30\u00f74=\", \"10\u00f76=\"),\n  @(\"32\u00f73=\", \"44\u00f75=\"),\n  @(\"56\u00f79=\", \"41\u00f72=\"),\n  @(\"89\u00f78=\", \"74\u00f78=\"),\n  @(\"26\u00f79=\", \"81\u00f78=\"),\n  @(\"60\u00f76=\", \"48\u00f75=\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($null, $null, $null, $null, $null, $null, $true, $null, $null, $null, 2) | Out-Null\n}\n"}
